$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(646, 163, 12, 227, 6, 8, 19, 0, 6)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $srcCell = $ws.Cells.Item($row, 4)
    $dstCell = $ws.Cells.Item($row, 5)
    $dstCell.Value = $values[$i]
    $srcCell.Copy()
    $dstCell.PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

$ws.Range("G9").Select()
